$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Hoja1")

# The "Estados de Cuenta" period table (rows 16-22, columns E/F/G) is
# refreshed: previous periods (2301..2307, ascending) are replaced by the
# new periods in descending order (2307..2301), each worker's "Valor Mora"
# carries along with its period (34666 stays with period 2307, the rest are
# 40000), and "Salario Basico" drops from 3000000 to 1000000 for every row.

$periods = @("2307", "2306", "2305", "2304", "2303", "2302", "2301")
$valorMora = @(34666, 40000, 40000, 40000, 40000, 40000, 40000)
$salarioBasico = @(1000000, 1000000, 1000000, 1000000, 1000000, 1000000, 1000000)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valorMora[$i]
    $ws.Range("G$row").Value = $salarioBasico[$i]
}
